# The deck ships two theme parts:
#   ppt/theme/theme1.xml  - currently "Office Theme" colours, wired only to
#                            the notes master relationship
#   ppt/theme/theme2.xml  - currently "Integral" colours, wired to the
#                            slide master / presentation (the theme that is
#                            actually visible on every slide)
#
# The authored change swaps the two parts' contents wholesale, so the
# slide master ends up using the plain "Office Theme" palette while the
# (functionally invisible) notes-master-only part ends up holding the old
# "Integral" palette. The <a:fontScheme> and <a:fmtScheme> blocks are
# byte-identical between the two theme parts already, so the only
# observable difference after the swap is the <a:clrScheme> (12 theme
# colours) used by the slide master / every slide.
#
# Re-point the live theme's colour scheme at the "Office Theme" palette
# using the standard PowerPoint theme-colour indices:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink

function ConvertTo-VbaRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeTheme = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

for ($i = 1; $i -le 12; $i++) {
    $tcs.Colors($i).RGB = ConvertTo-VbaRgb $officeTheme[$i - 1]
}
